$d = $word.ActiveDocument

$d.Content.Find.Execute("ia-2-inventory-whitelist.txt", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ia-2.1-inventory-whitelist.txt", 2)

$d.Content.Find.Execute("ia-2-admin-sudo.png", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ia-2.1-admin-sudo.png", 2)

$d.Content.Find.Execute("ia-4-two-factor-auth.png", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ia-2.1-two-factor-auth.png", 2)
